$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl25"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.835473666666667
$ws.Range("H2").Value = 8.506421
$ws.Range("I2").Value = 0.1177775131860914
$ws.Range("J2").Value = 0.1254491357467524
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.2443278446591134
$ws.Range("P2").Value = 0.3265937887468804
$ws.Range("Q2").Value = 1.617976093357556
$ws.Range("R2").Value = 14.561784840218
$ws.Range("S2").Value = 0.02877632594606803
$ws.Range("T2").Value = 0.04097090853855358

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl25"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.835473666666667
$ws.Range("H3").Value = 8.506421
$ws.Range("I3").Value = 0.1177775131860914
$ws.Range("J3").Value = 0.1254491357467524
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.7648465
$ws.Range("N3").Value = 3.529693
$ws.Range("O3").Value = 0.7556721553408867
$ws.Range("P3").Value = 0.6734062112531195
$ws.Range("Q3").Value = 5.004175776458833
$ws.Range("R3").Value = 30.025054658753
$ws.Range("S3").Value = 0.08900118724002343
$ws.Range("T3").Value = 0.08447822720819881

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl25"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.985602333333333
$ws.Range("H4").Value = 20.956807
$ws.Range("I4").Value = 0.2901620567311298
$ws.Range("J4").Value = 0.3090622161966227
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.5706193333333334
$ws.Range("N4").Value = 1.711858
$ws.Range("O4").Value = 0.2443278446591134
$ws.Range("P4").Value = 0.3265937887468804
$ws.Range("Q4").Value = 3.986119746378445
$ws.Range("R4").Value = 35.87507771740601
$ws.Range("S4").Value = 0.07089466992297232
$ws.Range("T4").Value = 0.1009378001461625

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl25"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.985602333333333
$ws.Range("H5").Value = 20.956807
$ws.Range("I5").Value = 0.2901620567311298
$ws.Range("J5").Value = 0.3090622161966227
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.7648465
$ws.Range("N5").Value = 3.529693
$ws.Range("O5").Value = 0.7556721553408867
$ws.Range("P5").Value = 0.6734062112531195
$ws.Range("Q5").Value = 12.32851582837517
$ws.Range("R5").Value = 73.971094970251
$ws.Range("S5").Value = 0.2192673868081575
$ws.Range("T5").Value = 0.2081244160504602

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ccl25"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.277146666666667
$ws.Range("H6").Value = 15.83144
$ws.Range("I6").Value = 0.2191976664868592
$ws.Range("J6").Value = 0.2334754493842435
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.5706193333333334
$ws.Range("N6").Value = 1.711858
$ws.Range("O6").Value = 0.2443278446591134
$ws.Range("P6").Value = 0.3265937887468804
$ws.Range("Q6").Value = 3.011241912835556
$ws.Range("R6").Value = 27.10117721552
$ws.Range("S6").Value = 0.05355609340704148
$ws.Range("T6").Value = 0.07625163159378061

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ccl25"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.277146666666667
$ws.Range("H7").Value = 15.83144
$ws.Range("I7").Value = 0.2191976664868592
$ws.Range("J7").Value = 0.2334754493842435
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.7648465
$ws.Range("N7").Value = 3.529693
$ws.Range("O7").Value = 0.7556721553408867
$ws.Range("P7").Value = 0.6734062112531195
$ws.Range("Q7").Value = 9.313353824653333
$ws.Range("R7").Value = 55.88012294792
$ws.Range("S7").Value = 0.1656415730798178
$ws.Range("T7").Value = 0.1572238177904629

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ccl25"
$ws.Range("C8").Value = "Ccr10"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.559845333333333
$ws.Range("H8").Value = 13.679536
$ws.Range("I8").Value = 0.1894030088117685
$ws.Range("J8").Value = 0.2017400700737227
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.5706193333333334
$ws.Range("N8").Value = 1.711858
$ws.Range("O8").Value = 0.2443278446591134
$ws.Range("P8").Value = 0.3265937887468804
$ws.Range("Q8").Value = 2.601935904209777
$ws.Range("R8").Value = 23.417423137888
$ws.Range("S8").Value = 0.04627642891493045
$ws.Range("T8").Value = 0.06588705382743826

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ccl25"
$ws.Range("C9").Value = "Ccr10"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.559845333333333
$ws.Range("H9").Value = 13.679536
$ws.Range("I9").Value = 0.1894030088117685
$ws.Range("J9").Value = 0.2017400700737227
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.7648465
$ws.Range("N9").Value = 3.529693
$ws.Range("O9").Value = 0.7556721553408867
$ws.Range("P9").Value = 0.6734062112531195
$ws.Range("Q9").Value = 8.047427077074666
$ws.Range("R9").Value = 48.28456246244799
$ws.Range("S9").Value = 0.1431265798968381
$ws.Range("T9").Value = 0.1358530162462845

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ccl25"
$ws.Range("C10").Value = "Ccr10"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.4167625
$ws.Range("H10").Value = 8.833525
$ws.Range("I10").Value = 0.183459754784151
$ws.Range("J10").Value = 0.1302731285986587
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.5706193333333334
$ws.Range("N10").Value = 1.711858
$ws.Range("O10").Value = 0.2443278446591134
$ws.Range("P10").Value = 0.3265937887468804
$ws.Range("Q10").Value = 2.520290073241667
$ws.Range("R10").Value = 15.12174043945
$ws.Range("S10").Value = 0.04482432646810108
$ws.Range("T10").Value = 0.04254639464094553

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Ccl25"
$ws.Range("C11").Value = "Ccr10"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.4167625
$ws.Range("H11").Value = 8.833525
$ws.Range("I11").Value = 0.183459754784151
$ws.Range("J11").Value = 0.1302731285986587
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.7648465
$ws.Range("N11").Value = 3.529693
$ws.Range("O11").Value = 0.7556721553408867
$ws.Range("P11").Value = 0.6734062112531195
$ws.Range("Q11").Value = 7.79490783945625
$ws.Range("R11").Value = 31.179631357825
$ws.Range("S11").Value = 0.13863542831605
$ws.Range("T11").Value = 0.08772673395771317

